# Natmi following Dr Hou advice
# Rebuilds the Il4-Il2rg LR-pairs sheet to include "M2" as an additional
# sending cluster (in addition to being a target cluster), and refreshes
# all of the derived statistics for the full 4x4 sending/target cluster
# matrix (ECs, FAPs, M2, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row in column order A..T:
#   A: Sending cluster       B: Ligand symbol            C: Receptor symbol
#   D: Target cluster        E: Ligand-expressing cells  F: Ligand detection rate
#   G: Ligand average expr   H: Ligand total expr        I: Ligand derived specificity (avg)
#   J: Ligand derived specificity (total)
#   K: Receptor-expressing cells   L: Receptor detection rate
#   M: Receptor average expr       N: Receptor total expr
#   O: Receptor derived specificity (avg)   P: Receptor derived specificity (total)
#   Q: Edge average expr weight    R: Edge total expr weight
#   S: Edge average expr derived specificity   T: Edge total expr derived specificity
$rowData = @(
    @("ECs", "Il4", "Il2rg", "ECs", 3, 1, 0.9334919999999999, 2.800476, 0.209095668602663, 0.209095668602663, 3, 1, 24.55401366666667, 73.662041, 0.2501033965205259, 0.250103396520526, 22.920975325724, 206.288777931516, 0.05229553691525631, 0.05229553691525632),
    @("ECs", "Il4", "Il2rg", "FAPs", 3, 1, 0.9334919999999999, 2.800476, 0.209095668602663, 0.209095668602663, 3, 1, 1.205204666666667, 3.615614, 0.01227602886956614, 0.01227602886956614, 1.125048914696, 10.125440232264, 0.002566864464267526, 0.002566864464267526),
    @("ECs", "Il4", "Il2rg", "M2", 3, 1, 0.9334919999999999, 2.800476, 0.209095668602663, 0.209095668602663, 3, 1, 67.337282, 202.011846, 0.6858871697837075, 0.6858871697837076, 62.85881404874399, 565.7293264386959, 0.1434160363519126, 0.1434160363519126),
    @("ECs", "Il4", "Il2rg", "sCs", 3, 1, 0.9334919999999999, 2.800476, 0.209095668602663, 0.209095668602663, 3, 1, 5.078950333333333, 15.236851, 0.0517334048262004, 0.05173340482620041, 4.741159504563999, 42.670435541076, 0.01081723087122661, 0.01081723087122661),
    @("FAPs", "Il4", "Il2rg", "ECs", 3, 1, 1.964334, 5.893002, 0.4399970552387631, 0.4399970552387631, 3, 1, 24.55401366666667, 73.662041, 0.2501033965205259, 0.250103396520526, 48.232283881898, 434.090554937082, 0.1100447579742441, 0.1100447579742441),
    @("FAPs", "Il4", "Il2rg", "FAPs", 3, 1, 1.964334, 5.893002, 0.4399970552387631, 0.4399970552387631, 3, 1, 1.205204666666667, 3.615614, 0.01227602886956614, 0.01227602886956614, 2.367424503692, 21.306820533228, 0.005401416552635145, 0.005401416552635145),
    @("FAPs", "Il4", "Il2rg", "M2", 3, 1, 1.964334, 5.893002, 0.4399970552387631, 0.4399970552387631, 3, 1, 67.337282, 202.011846, 0.6858871697837075, 0.6858871697837076, 132.272912500188, 1190.456212501692, 0.3017883349308809, 0.3017883349308809),
    @("FAPs", "Il4", "Il2rg", "sCs", 3, 1, 1.964334, 5.893002, 0.4399970552387631, 0.4399970552387631, 3, 1, 5.078950333333333, 15.236851, 0.0517334048262004, 0.05173340482620041, 9.976754824078, 89.790793416702, 0.02276254578100299, 0.02276254578100299),
    @("M2", "Il4", "Il2rg", "ECs", 2, 0.6666666666666666, 0.1531806666666667, 0.459542, 0.03431139625585257, 0.03431139625585256, 3, 1, 24.55401366666667, 73.662041, 0.2501033965205259, 0.250103396520526, 3.761200182802444, 33.850801645222, 0.008581396742950383, 0.008581396742950383),
    @("M2", "Il4", "Il2rg", "FAPs", 2, 0.6666666666666666, 0.1531806666666667, 0.459542, 0.03431139625585257, 0.03431139625585256, 3, 1, 1.205204666666667, 3.615614, 0.01227602886956614, 0.01227602886956614, 0.1846140543097778, 1.661526488788, 0.0004212076909919698, 0.0004212076909919697),
    @("M2", "Il4", "Il2rg", "M2", 2, 0.6666666666666666, 0.1531806666666667, 0.459542, 0.03431139625585257, 0.03431139625585256, 3, 1, 67.337282, 202.011846, 0.6858871697837075, 0.6858871697837076, 10.31476974828133, 92.83292773453199, 0.02353374646925402, 0.02353374646925402),
    @("M2", "Il4", "Il2rg", "sCs", 2, 0.6666666666666666, 0.1531806666666667, 0.459542, 0.03431139625585257, 0.03431139625585256, 3, 1, 5.078950333333333, 15.236851, 0.0517334048262004, 0.05173340482620041, 0.7779969980268888, 7.001972982242, 0.001775045352656198, 0.001775045352656198),
    @("sCs", "Il4", "Il2rg", "ECs", 3, 1, 1.413418666666667, 4.240256, 0.3165958799027214, 0.3165958799027214, 3, 1, 24.55401366666667, 73.662041, 0.2501033965205259, 0.250103396520526, 34.70510125805512, 312.345911322496, 0.07918170488807515, 0.07918170488807515),
    @("sCs", "Il4", "Il2rg", "FAPs", 3, 1, 1.413418666666667, 4.240256, 0.3165958799027214, 0.3165958799027214, 3, 1, 1.205204666666667, 3.615614, 0.01227602886956614, 0.01227602886956614, 1.703458773020445, 15.331128957184, 0.003886540161671503, 0.003886540161671503),
    @("sCs", "Il4", "Il2rg", "M2", 3, 1, 1.413418666666667, 4.240256, 0.3165958799027214, 0.3165958799027214, 3, 1, 67.337282, 202.011846, 0.6858871697837075, 0.6858871697837076, 95.17577134139735, 856.5819420725761, 0.2171490520316601, 0.2171490520316601),
    @("sCs", "Il4", "Il2rg", "sCs", 3, 1, 1.413418666666667, 4.240256, 0.3165958799027214, 0.3165958799027214, 3, 1, 5.078950333333333, 15.236851, 0.0517334048262004, 0.05173340482620041, 7.178683208206223, 64.60814887385601, 0.01637858282131461, 0.01637858282131461)
)

$r = 2
foreach ($row in $rowData) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
